# Refined metadata to be additional tab
#
# 1. Refresh the "time_taken" column (F2:F35) on the existing "data" sheet
#    with the new query timestamps.
# 2. Add a new "metadata" worksheet (after "data") describing the panel
#    query itself (name, id, version, version timestamp, query time, and
#    the API request URL used to fetch it).

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# --- 1. Updated time_taken values for rows 2..35 -------------------------
$timestamps = @(
    "2021-10-05 14:33:08.783606",
    "2021-10-05 14:33:08.783614",
    "2021-10-05 14:33:08.783617",
    "2021-10-05 14:33:08.783620",
    "2021-10-05 14:33:08.783623",
    "2021-10-05 14:33:08.783626",
    "2021-10-05 14:33:08.783628",
    "2021-10-05 14:33:08.783631",
    "2021-10-05 14:33:08.783634",
    "2021-10-05 14:33:08.783636",
    "2021-10-05 14:33:08.783639",
    "2021-10-05 14:33:08.783642",
    "2021-10-05 14:33:08.783644",
    "2021-10-05 14:33:08.783647",
    "2021-10-05 14:33:08.783650",
    "2021-10-05 14:33:08.783652",
    "2021-10-05 14:33:08.783655",
    "2021-10-05 14:33:08.783658",
    "2021-10-05 14:33:08.783660",
    "2021-10-05 14:33:08.783663",
    "2021-10-05 14:33:08.783665",
    "2021-10-05 14:33:08.783668",
    "2021-10-05 14:33:08.783670",
    "2021-10-05 14:33:08.783673",
    "2021-10-05 14:33:08.783676",
    "2021-10-05 14:33:08.783679",
    "2021-10-05 14:33:08.783681",
    "2021-10-05 14:33:08.783684",
    "2021-10-05 14:33:08.783686",
    "2021-10-05 14:33:08.783689",
    "2021-10-05 14:33:08.783691",
    "2021-10-05 14:33:08.783694",
    "2021-10-05 14:33:08.783697",
    "2021-10-05 14:33:08.783699"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $data.Range("F$row").Value = $timestamps[$i]
}

# --- 2. New "metadata" sheet, placed after "data" -------------------------
$meta = $wb.Worksheets.Add([Type]::Missing, $data)
$meta.Name = "metadata"

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Angelman Rett like syndromes"
$meta.Range("C2").Value = 41
# Force "1.0" to stay text (otherwise Excel auto-coerces it to the number 1).
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.0"
$meta.Range("E2").Value = "2021-06-07T00:34:46.430247Z"
$meta.Range("F2").Value = "2021-10-05 14:33:08.779680"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/41/?format=json"

# Reuse the existing bold/bordered header style from the "data" sheet (its
# header row - and the leading index column - use the same formatting) so
# no new style entries are created.
$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Keep "data" as the active/selected sheet (only a tab was added).
$data.Activate()
